$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.502.87'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '3.977.08'
$ws.Range("E3").Value = '  -2.05%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'541.78"
$ws.Range("E5").Value = '  +4.41%  '
$ws.Range("D6").Value = "'149.05"
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("D7").Value = '3.969.34'
$ws.Range("E7").Value = '  -2.04%  '
$ws.Range("D8").Value = "'0.688"
$ws.Range("E8").Value = '  -6.44%  '
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  -4.06%  '
$ws.Range("E11").Value = '  -4.92%  '
$ws.Range("D12").Value = "'56.49"
$ws.Range("E12").Value = '  +18.54%  '
$ws.Range("D13").Value = "'0.0000320"
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("E14").Value = '  -3.71%  '
$ws.Range("D15").Value = '4.606.01'
$ws.Range("E15").Value = '  -2.10%  '
$ws.Range("D16").Value = '3.962.53'
$ws.Range("E16").Value = '  -2.58%  '
$ws.Range("D17").Value = "'13.94"
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("D18").Value = "'20.65"
$ws.Range("E18").Value = '  -2.60%  '
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("E20").Value = '  -2.63%  '
$ws.Range("D21").Value = '71.332.60'
$ws.Range("E21").Value = '  -1.69%  '
$ws.Range("D22").Value = "'428.03"
$ws.Range("E22").Value = '  -3.62%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = "'97.68"
$ws.Range("E23").Value = '  -6.92%  '
$ws.Range("B24").Value = 'ImmutableX'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D24").Value = "'3.59"
$ws.Range("E24").Value = '  +0.86%  '
$ws.Range("D25").Value = "'4.23"
$ws.Range("E25").Value = '  +5.62%  '
$ws.Range("D26").Value = "'14.46"
$ws.Range("E26").Value = '  -2.24%  '
$ws.Range("D27").Value = "'11.50"
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("E28").Value = '  -2.34%  '
$ws.Range("D29").Value = "'3.73"
$ws.Range("E29").Value = '  +13.81%  '
$ws.Range("D30").Value = "'5.90"
$ws.Range("E30").Value = '  +1.46%  '
$ws.Range("D31").Value = "'36.69"
$ws.Range("E31").Value = '  -3.13%  '
$ws.Range("D32").Value = "'7.74"
$ws.Range("E32").Value = '  +13.79%  '
$ws.Range("D33").Value = "'51.42"
$ws.Range("E33").Value = '  +20.95%  '
$ws.Range("D34").Value = "'695.08"
$ws.Range("E34").Value = '  +1.93%  '
$ws.Range("B35").Value = 'Cosmos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D35").Value = "'13.45"
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = "'0.131"
$ws.Range("E36").Value = '  +0.85%  '
$ws.Range("D37").Value = "'65.56"
$ws.Range("E37").Value = '  -2.25%  '
$ws.Range("E38").Value = '  +2.54%  '
$ws.Range("D39").Value = '0.0₃0830'
$ws.Range("E39").Value = '  -3.62%  '
$ws.Range("E40").Value = '  +0.57%  '
$ws.Range("E41").Value = '  -2.39%  '
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("D44").Value = "'3.28"
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("D45").Value = "'0.0486"
$ws.Range("E45").Value = '  -2.36%  '
$ws.Range("E46").Value = '  -5.69%  '
$ws.Range("E47").Value = '  +1.07%  '
$ws.Range("D48").Value = "'9.77"
$ws.Range("E48").Value = '  +6.35%  '
$ws.Range("E49").Value = '  -3.50%  '
$ws.Range("D50").Value = "'3.01"
$ws.Range("E50").Value = '  -1.88%  '
$ws.Range("E51").Value = '  +1.87%  '

# Reset style on quote-prefixed cells so no stray number-format/style is introduced
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
